# Auto-generated edit script: updates crypto price/volume table values
# per commit "Updated cryptos list on Thu Aug 17 23:42:53 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.765.46'
$ws.Range('E2').Value = '  -7.24%  '
$ws.Range('D3').Value = '1.695.64'
$ws.Range('E3').Value = '  -6.45%  '
$ws.Range('D4').Value = "'1.005"
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'220.56"
$ws.Range('E5').Value = '  -5.07%  '
$ws.Range('D6').Value = "'0.5118"
$ws.Range('E6').Value = '  -13.35%  '
$ws.Range('D7').Value = "'1.006"
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = "'0.2576"
$ws.Range('E8').Value = '  -6.44%  '
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').Value = "'0.06175"
$ws.Range('E10').Value = '  -8.38%  '
$ws.Range('D11').Value = "'0.07348"
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').Value = '1.701.96'
$ws.Range('E12').Value = '  -6.22%  '
$ws.Range('D13').Value = "'4.463"
$ws.Range('E13').Value = '  -4.43%  '
$ws.Range('D14').Value = "'0.5802"
$ws.Range('E14').Value = '  -7.01%  '
$ws.Range('D15').Value = '1.921.73'
$ws.Range('E15').Value = '  -6.96%  '
$ws.Range('D16').Value = "'0.000008168"
$ws.Range('E16').Value = '  -12.27%  '
$ws.Range('D17').Value = "'65.25"
$ws.Range('E17').Value = '  -12.49%  '
$ws.Range('D18').Value = '26.768.50'
$ws.Range('E18').Value = '  -6.45%  '
$ws.Range('D19').Value = "'4.991"
$ws.Range('E19').Value = '  -8.25%  '
$ws.Range('D20').Value = "'1.005"
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -6.28%  '
$ws.Range('D22').Value = "'186.64"
$ws.Range('E22').Value = '  -10.43%  '
$ws.Range('D23').Value = "'6.265"
$ws.Range('E23').Value = '  -7.38%  '
$ws.Range('D24').Value = "'1.006"
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = "'142.30"
$ws.Range('E25').Value = '  -7.76%  '
$ws.Range('D26').Value = "'7.497"
$ws.Range('E26').Value = '  -3.80%  '
$ws.Range('D27').Value = "'0.1146"
$ws.Range('E27').Value = '  -9.81%  '
$ws.Range('D28').Value = "'15.19"
$ws.Range('E28').Value = '  -6.81%  '
$ws.Range('D29').Value = "'1.332"
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -6.61%  '
$ws.Range('D31').Value = "'1.348"
$ws.Range('E31').Value = '  -5.60%  '
$ws.Range('D32').Value = "'3.448"
$ws.Range('E32').Value = '  -7.56%  '
$ws.Range('D33').Value = "'3.428"
$ws.Range('E33').Value = '  -7.03%  '
$ws.Range('D34').Value = "'1.638"
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').Value = "'0.9880"
$ws.Range('E35').Value = '  -5.92%  '
$ws.Range('D36').Value = "'2.420"
$ws.Range('E36').Value = '  -4.00%  '
$ws.Range('D37').Value = "'0.5966"
$ws.Range('E37').Value = '  -6.02%  '
$ws.Range('D38').Value = "'2.661"
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.089.44'
$ws.Range('E39').Value = '  -3.72%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.01593"
$ws.Range('E40').Value = '  -5.64%  '
$ws.Range('D41').Value = "'0.8609"
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').Value = "'5.850"
$ws.Range('E42').Value = '  -9.01%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').Value = "'97.56"
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').Value = '1.845.26'
$ws.Range('E45').Value = '  -6.30%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000106"
$ws.Range('E46').Value = '  -5.25%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'55.94"
$ws.Range('E47').Value = '  -7.50%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = "'1.004"
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  -3.94%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = "'0.4308"
$ws.Range('E50').Value = '  -4.54%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'7.915"
$ws.Range('E51').Value = '  -4.23%  '
